$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (first sample record): rename Joe Doe -> Jony L, and USA -> United States
$ws.Range("B2").Value = "JonyL"
$ws.Range("J2").Value = "Jony"
$ws.Range("K2").Value = "L"
$ws.Range("O2").Value = "United States"

# Row 3 (second sample record): rename Harry Smith -> Shyam P, update DOB year,
# and change country from USA to India
$ws.Range("B3").Value = "ShyamP"
$ws.Range("G3").Value = 2000
$ws.Range("J3").Value = "Shyam"
$ws.Range("K3").Value = "P"
$ws.Range("O3").Value = "India"
